$wb = $excel.ActiveWorkbook

# --- "Summary B" sheet (xl/worksheets/sheet2.xml): add new Mean/SD summary rows ---
$ws = $wb.Worksheets.Item("Summary B")

# Row 52 first so "Ag Mean" becomes shared-string index 51 ...
$ws.Cells.Item(52,1).Value = "Ag Mean"
$ws.Cells.Item(52,3).Formula = "=AVERAGE(C5,C9,C17,C21,C25,C29,C33,C37,C41,C45)"
$ws.Cells.Item(52,4).Formula = "=AVERAGE(D5,D9,D17,D21,D25,D29,D33,D37,D41,D45)"
$ws.Cells.Item(52,5).Formula = "=AVERAGE(E5,E9,E17,E21,E25,E29,E33,E37,E41,E45)"
$ws.Cells.Item(52,6).Formula = "=AVERAGE(F5,F9,F17,F21,F25,F29,F33,F37,F41,F45)"
$ws.Cells.Item(52,7).Formula = "=AVERAGE(G5,G9,G13,G17,G21,G25,G29,G33,G37,G41)"
$ws.Cells.Item(52,8).Formula = "=AVERAGE(H5,H9,H13,H17,H21,H25,H29,H37,H41,H45)"
$ws.Cells.Item(52,9).Formula = "=AVERAGE(I5,I9,I13,I17,I21,I25,I29,I37,I41,I45)"
$ws.Cells.Item(52,10).Formula = "=AVERAGE(J5,J9,J13,J17,J21,J37,J41,J45)"
$ws.Cells.Item(52,11).Formula = "=AVERAGE(K5,K9,K13,K17,K21,K25,K29,K33,K37,K41,K45)"
$ws.Cells.Item(52,12).Formula = "=AVERAGE(L5,L9,L13,L17,L21,L25,L29,L33,L37,L41,L45)"
$ws.Cells.Item(52,13).Formula = "=AVERAGE(M5,M9,M13,M17,M21,M25,M29,M33,M37,M41,M45)"
$ws.Cells.Item(52,14).Formula = "=AVERAGE(N5,N9,N13,N17,N21,N25,N29,N33,N37,N41,N45)"
$ws.Cells.Item(52,15).Formula = "=AVERAGE(O5,O9,O13,O17,O21,O25,O29,O33,O37,O41,O45)"

# ... then row 50 so "Prairie Mean" becomes shared-string index 52
$ws.Cells.Item(50,1).Value = "Prairie Mean"
$ws.Cells.Item(50,3).Formula = "=AVERAGE(C3,C7,C15,C19,C23,C27,C31,C35,C39,C43)"
$ws.Cells.Item(50,4).Formula = "=AVERAGE(D3,D7,D15,D19,D23,D27,D31,D35,D39,D43)"
$ws.Cells.Item(50,5).Formula = "=AVERAGE(E3,E7,E15,E19,E23,E27,E31,E35,E39,E43)"
$ws.Cells.Item(50,6).Formula = "=AVERAGE(F3,F7,F15,F19,F23,F27,F31,F35,F39,F43)"
$ws.Cells.Item(50,7).Formula = "=AVERAGE(G3,G7,G11,G15,G19,G23,G27,G31,G35,G39)"
$ws.Cells.Item(50,8).Formula = "=AVERAGE(H3,H7,H11,H15,H19,H23,H27,H35,H39,H43)"
$ws.Cells.Item(50,9).Formula = "=AVERAGE(I3,I7,I11,I15,I19,I23,I27,I35,I39,I43)"
$ws.Cells.Item(50,10).Formula = "=AVERAGE(J3,J7,J11,J15,J19,J35,J39,J43)"
$ws.Cells.Item(50,11).Formula = "=AVERAGE(K3,K7,K11,K15,K19,K23,K27,K31,K35,K39,K43)"
$ws.Cells.Item(50,12).Formula = "=AVERAGE(L3,L7,L11,L15,L19,L23,L27,L31,L35,L39,L43)"
$ws.Cells.Item(50,13).Formula = "=AVERAGE(M3,M7,M11,M15,M19,M23,M27,M31,M35,M39,M43)"
$ws.Cells.Item(50,14).Formula = "=AVERAGE(N3,N7,N11,N15,N19,N23,N27,N31,N35,N39,N43)"
$ws.Cells.Item(50,15).Formula = "=AVERAGE(O3,O7,O11,O15,O19,O23,O27,O31,O35,O39,O43)"

# Row 51 - "Prairie SD" (re-uses existing shared string)
$ws.Cells.Item(51,1).Value = "Prairie SD"
$ws.Cells.Item(51,3).Formula = "=AVERAGE(C4,C8,C16,C20,C24,C28,C32,C36,C40,C44)"
$ws.Cells.Item(51,4).Formula = "=AVERAGE(D4,D8,D16,D20,D24,D28,D32,D36,D40,D44)"
$ws.Cells.Item(51,5).Formula = "=AVERAGE(E4,E8,E16,E20,E24,E28,E32,E36,E40,E44)"
$ws.Cells.Item(51,6).Formula = "=AVERAGE(F4,F8,F16,F20,F24,F28,F32,F36,F40,F44)"
$ws.Cells.Item(51,7).Formula = "=AVERAGE(G4,G8,G12,G16,G20,G24,G28,G32,G36,G40)"
$ws.Cells.Item(51,8).Formula = "=AVERAGE(H4,H8,H12,H16,H20,H24,H28,H36,H40,H44)"
$ws.Cells.Item(51,9).Formula = "=AVERAGE(I4,I8,I12,I16,I20,I24,I28,I36,I40,I44)"
$ws.Cells.Item(51,10).Formula = "=AVERAGE(J4,J8,J12,J16,J20)"
$ws.Cells.Item(51,11).Formula = "=AVERAGE(K4,K8,K12,K16,K20,K24,K28,K32,K36,K40,K44)"
$ws.Cells.Item(51,12).Formula = "=AVERAGE(L4,L8,L12,L16,L20,L24,L28,L32,L36,L40,L44)"
$ws.Cells.Item(51,13).Formula = "=AVERAGE(M4,M8,M12,M16,M20,M24,M28,M32,M36,M40,M44)"
$ws.Cells.Item(51,14).Formula = "=AVERAGE(N4,N8,N12,N16,N20,N24,N28,N32,N36,N40,N44)"
$ws.Cells.Item(51,15).Formula = "=AVERAGE(O4,O8,O12,O16,O20,O24,O28,O32,O36,O40,O44)"

# Row 53 - "Ag SD" (re-uses existing shared string)
$ws.Cells.Item(53,1).Value = "Ag SD"
$ws.Cells.Item(53,3).Formula = "=AVERAGE(C6,C10,C18,C22,C26,C30,C34,C38,C42,C46)"
$ws.Cells.Item(53,4).Formula = "=AVERAGE(D6,D10,D18,D22,D26,D30,D34,D38,D42,D46)"
$ws.Cells.Item(53,5).Formula = "=AVERAGE(E6,E10,E18,E22,E26,E30,E34,E38,E42,E46)"
$ws.Cells.Item(53,6).Formula = "=AVERAGE(F6,F10,F18,F22,F26,F30,F34,F38,F42,F46)"
$ws.Cells.Item(53,7).Formula = "=AVERAGE(G6,G10,G14,G18,G22,G26,G30,G34,G38,G42)"
$ws.Cells.Item(53,8).Formula = "=AVERAGE(H6,H10,H14,H18,H22,H26,H30,H38,H42,H46)"
$ws.Cells.Item(53,9).Formula = "=AVERAGE(I6,I10,I14,I18,I22,I26,I30,I38,I42,I46)"
$ws.Cells.Item(53,10).Formula = "=AVERAGE(J6,J10,J14,J18,J22,J38,J42,J46)"
$ws.Cells.Item(53,11).Formula = "=AVERAGE(K6,K10,K14,K18,K22,K26,K30,K34,K38,K42,K46)"
$ws.Cells.Item(53,12).Formula = "=AVERAGE(L6,L10,L14,L18,L22,L26,L30,L34,L38,L42,L46)"
$ws.Cells.Item(53,13).Formula = "=AVERAGE(M6,M10,M14,M18,M22,M26,M30,M34,M38,M42,M46)"
$ws.Cells.Item(53,14).Formula = "=AVERAGE(N6,N10,N14,N18,N22,N26,N30,N34,N38,N42,N46)"
$ws.Cells.Item(53,15).Formula = "=AVERAGE(O6,O10,O14,O18,O22,O26,O30,O34,O38,O42,O46)"

# "Summary B" becomes the active tab, scrolled near the new rows, selection on J51
$ws.Activate()
$ws.Range("J51").Select()

# --- "Summary" sheet (xl/worksheets/sheet3.xml): selection moves off of F14 to H24 ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("H24").Select()

# Restore "Summary B" as the active/selected sheet+cell (mirrors final saved view state)
$ws.Activate()
$ws.Range("J51").Select()
